$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 15 entirely (duplicate "Austin-Power-BI-User-Group" row) - shifts rows 16-21 up by one.
$ws.Rows(15).Delete()

# Add "Errors" to G6 (shared string index 2 => "Errors")
$ws.Range("G6").Value = "Errors"

Write-Output "Done"
